$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 89 (existing rows 89:122 shift down to 95:128)
$ws.Rows("89:94").Insert()

# Common / repeated column values for this product block
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad = "Cultivar IV Región"

# New weekly data rows (fecha serial 44511) for rows 89-94
$newRows = @(
    @{ Row = 89; Calidad = "Cuarta";                   Volumen = 250; PMin = 1200;  PMax = 1200;  PProm = 1200;  Unidad = "`$/kilo (en caja de 15 kilos)"; Origen = "Provincia de Limarí"; PrecioKg = 1200; KgUnidad = 1 },
    @{ Row = 90; Calidad = "Especial";                 Volumen = 280; PMin = 20000; PMax = 20000; PProm = 20000; Unidad = "`$/bandeja 8 kilos";             Origen = "Provincia de Limarí"; PrecioKg = 2500; KgUnidad = 8 },
    @{ Row = 91; Calidad = "Extra (doble especial)";   Volumen = 350; PMin = 24000; PMax = 24000; PProm = 24000; Unidad = "`$/bandeja 8 kilos";             Origen = "Provincia de Limarí"; PrecioKg = 3000; KgUnidad = 8 },
    @{ Row = 92; Calidad = "Primera";                  Volumen = 330; PMin = 16000; PMax = 16000; PProm = 16000; Unidad = "`$/bandeja 8 kilos";             Origen = "Provincia de Limarí"; PrecioKg = 2000; KgUnidad = 8 },
    @{ Row = 93; Calidad = "Segunda";                  Volumen = 300; PMin = 14400; PMax = 14400; PProm = 14400; Unidad = "`$/bandeja 8 kilos";             Origen = "Provincia de Limarí"; PrecioKg = 1800; KgUnidad = 8 },
    @{ Row = 94; Calidad = "Tercera";                  Volumen = 220; PMin = 1400;  PMax = 1400;  PProm = 1400;  Unidad = "`$/kilo (en caja de 15 kilos)"; Origen = "Provincia de Limarí"; PrecioKg = 1400; KgUnidad = 1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = 44511
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
